$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Characters")
$ws2 = $wb.Worksheets.Item("Skills")
$ws5 = $wb.Worksheets.Item("Buff")

# ---------------------------------------------------------------------------
# New shared strings must be created in this exact order so the rebuilt
# sharedStrings table lines up with the target workbook (unused strings are
# garbage-collected and new ones are appended in write order).
# ---------------------------------------------------------------------------

# 1) Buff sheet: new "bleed" debuff row (row 6)
$ws5.Range("A6").Value = "출혈"
$ws5.Range("B6").Value = "출혈 스킬에 공격당할 시 지속데미지를 받는다."

# 2) Skills sheet: new castingTime column header
$ws2.Range("J1").Value = "castingTime"

# 3) Skills sheet: rogue skill renamed from "헤이스트" to "암살"
$ws2.Range("B8").Value = "암살"

# 4) Buff sheet: new column headers
$ws5.Range("F1").Value = "durationTime"
$ws5.Range("G1").Value = "trueDamage"

# 5) Skills sheet: new animTime column header
$ws2.Range("I1").Value = "animTime"

# ---------------------------------------------------------------------------
# Remaining data edits (reuse existing shared strings / plain numbers).
# ---------------------------------------------------------------------------

# The skill that used to be row 8 ("헤이스트") is now row 9; restore its stats.
$ws2.Range("B9").Value = "헤이스트"
$ws2.Range("E9").Value = 0
$ws2.Range("F9").Value = 0

# Row 8 ("암살", formerly "일섬" data) takes on the old row-9 stats.
$ws2.Range("E8").Value = 10
$ws2.Range("F8").Value = 1.2

# New animTime / castingTime values for every skill row.
$ws2.Range("I2").Value = 1
$ws2.Range("J2").Value = 100
$ws2.Range("I3").Value = 1
$ws2.Range("J3").Value = 100
$ws2.Range("I4").Value = 1
$ws2.Range("J4").Value = 100
$ws2.Range("I5").Value = 1
$ws2.Range("J5").Value = 100
$ws2.Range("I6").Value = 1
$ws2.Range("J6").Value = 100
$ws2.Range("I7").Value = 1
$ws2.Range("J7").Value = 100
$ws2.Range("I8").Value = 1
$ws2.Range("J8").Value = 100
$ws2.Range("I9").Value = 1
$ws2.Range("J9").Value = 100
$ws2.Range("I10").Value = 1
$ws2.Range("J10").Value = 100

# New durationTime / trueDamage values for the existing buffs.
$ws5.Range("F2").Value = 60
$ws5.Range("G2").Value = 0
$ws5.Range("F3").Value = 30
$ws5.Range("G3").Value = 0
$ws5.Range("F4").Value = 30
$ws5.Range("G4").Value = 0
$ws5.Range("F5").Value = 30
$ws5.Range("G5").Value = 0

# Finish the new "출혈" row's numeric columns.
$ws5.Range("C6").Value = 0
$ws5.Range("D6").Value = 0
$ws5.Range("E6").Value = 0
$ws5.Range("F6").Value = 3
$ws5.Range("G6").Value = 1

# Widen the Buff sheet's durationTime column a bit (best-effort; COM column
# widths are stored in character units and round to the nearest pixel).
$ws5.Columns.Item(6).ColumnWidth = 11.7142857142857

# ---------------------------------------------------------------------------
# View / selection state.
# ---------------------------------------------------------------------------

$ws2.Activate()
$ws2.Range("M8").Select()

$ws5.Activate()
$ws5.Range("G2").Select()

$ws1.Activate()
$ws1.Range("K10").Select()
